$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "N" column header in the totals row to show the overall N
$ws.Range("R77").Value = "Nges = 1642"

# Add a new row (78) with each cluster's share of the overall N, plus a
# 100% check total in column R
$ws.Range("A78").Formula = "=A77/1642"
$ws.Range("B78:Q78").Formula = "=B77/1642"
$ws.Range("R78").Value = 1

$ws.Range("R78").NumberFormat = "0%"
$ws.Range("A78:Q78").NumberFormat = "0.000%"

# Scroll the view down to where the action now is
$ws.Application.Goto $ws.Range("H87"), $false
$ws.Range("H87").Select()
$excel.ActiveWindow.ScrollRow = 42
